# Auto-generated Excel COM-interop script
# Applies crypto price/volume/name/link updates per the commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.286.24'
$ws.Range('E2').Value = '  +1.87%  '
$ws.Range('D3').Value = '1.875.47'
$ws.Range('E3').Value = '  +4.16%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.93'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.80%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.06%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5072'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.97%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3939'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.21%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.09642'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.77%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.146'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.82%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '40.94'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.12%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.498'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.44%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.04'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.01%  '
$ws.Range('D14').Value = '1.877.94'
$ws.Range('E14').Value = '  +4.55%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.438'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.58%  '
$ws.Range('B16').Value = 'BinanceUSD'
$ws.Range('C16').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.001'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.11%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001134'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.40%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '93.01'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.80%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06610'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.75%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.66'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.44%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.001'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.08%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.195'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.87%  '
$ws.Range('D23').Value = '28.338.58'
$ws.Range('E23').Value = '  +1.95%  '
$ws.Range('E24').Value = '  +2.95%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.303'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.34%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.574'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +7.19%  '
$ws.Range('D27').Value = '2.094.64'
$ws.Range('E27').Value = '  +4.56%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '21.28'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.44%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '158.88'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.69%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '127.69'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.97%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.1069'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.07%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.071'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.67%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.641'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.16%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.622'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.51%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.571'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +7.92%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06746'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.80%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02401'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.72%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2194'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.91%  '
$ws.Range('B39').Value = 'Aptos'
$ws.Range('C39').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '11.53'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.73%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6383'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.83%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.009'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.185'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.22%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.001'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.50'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.80%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6002'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.56%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.659'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.265'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.31%  '
$ws.Range('E48').Value = '  +3.65%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '124.28'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.82%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.197'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.32%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06854'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.21%  '
